# Update calibration data on Blad1 (measured "aflezing" readings in column B
# and the derived "Stroom" values in column C), then move the active
# selection on both sheets to match the author's last-saved cursor position.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Blad1")
$ws2 = $wb.Worksheets.Item("Blad2")

# row -> (B value, C value)
$data = @{
    2  = @(13,   2.0169999999999999)
    3  = @(27,   2.1059999999999999)
    4  = @(56,   2.3159999999999998)
    5  = @(145,  2.58)
    6  = @(296,  2.8319999999999999)
    7  = @(450,  2.992)
    8  = @(608,  3.1320000000000001)
    9  = @(764,  3.2269999999999999)
    10 = @(926,  3.319)
    11 = @(1086, 3.4020000000000001)
    12 = @(1251, 3.476)
    13 = @(1412, 3.54)
    14 = @(1568, 3.6)
    15 = @(1740, 3.66)
    16 = @(1905, 3.7130000000000001)
    17 = @(2068, 3.7639999999999998)
    18 = @(2240, 3.8140000000000001)
    19 = @(2398, 3.8580000000000001)
    20 = @(2571, 3.903)
    21 = @(2731, 3.944)
    22 = @(2906, 3.988)
    23 = @(3071, 4.0259999999999998)
    24 = @(3246, 4.0659999999999998)
    25 = @(3416, 4.1020000000000003)
    26 = @(3583, 4.1369999999999996)
    27 = @(3755, 4.1769999999999996)
    28 = @(3917, 4.2030000000000003)
    29 = @(4097, 4.2350000000000003)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws1.Cells.Item($row, 2).Value = $vals[0]
    $ws1.Cells.Item($row, 3).Value = $vals[1]
}

# Restore the sheets' saved cursor positions without changing which sheet
# is the currently active (tab-selected) one -- Blad1 stays active.
$ws2.Range("B1").Select()
$ws1.Select()
$ws1.Range("D20").Select()
